# "with all federal holidays added!"
#
# The "1-15" sheet is a July-2021 CBOC sign-in grid. Sunday July 4th
# (columns H/I) was already marked as a holiday (gray fill + "X").
# Independence Day was observed on Monday July 5th (columns J/K), so
# this change extends the same holiday styling/marker to the Monday
# column pair for every site row, and narrows columns J/K to match the
# other weekend/holiday columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1-15")

# Narrow columns J and K from 4.5 to 2.5 characters (matches column I's
# width, the "holiday" look of the grid). 1.67 in the COM ColumnWidth
# units corresponds to a stored xlsx width of 2.5 for this workbook.
$ws.Range("J1").ColumnWidth = 1.67
$ws.Range("K1").ColumnWidth = 1.67

# Rows that get the "X" marker written into the newly-highlighted
# Monday columns (every row that already shows "X" under Sunday).
$xRows = @(5,6,8,9,11,12,14,15,17,18,20,21,23,24,26,27)

for ($r = 2; $r -le 27; $r++) {
    # Column H/I already carry the correct holiday formatting (gray
    # fill matching the weekend columns) for this row; clone it onto
    # J/K without touching whatever value J/K already hold.
    $ws.Range("H$r").Copy()
    $ws.Range("J$r").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
    $ws.Range("I$r").Copy()
    $ws.Range("K$r").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
    $excel.CutCopyMode = $false

    if ($xRows -contains $r) {
        $ws.Range("J$r").Value = "X"
        $ws.Range("K$r").Value = "X"
    }
}
